# Applies the per-cell Price (D) and Volume(1h) (E) updates described in the commit diff.
# A leading apostrophe forces Excel to store the value as text (matching the original
# inline-string cell type) instead of auto-converting numeric-looking strings to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'70.843.45"
$ws.Range('E2').Value = "'  +2.23%  "
$ws.Range('D3').Value = "'3.548.99"
$ws.Range('E3').Value = "'  +1.01%  "
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = "'  -0.03%  "
$ws.Range('D5').Value = "'607.63"
$ws.Range('E5').Value = "'  +4.25%  "
$ws.Range('D6').Value = "'174.87"
$ws.Range('E6').Value = "'  +0.54%  "
$ws.Range('E7').Value = "'  -0.55%  "
$ws.Range('D8').Value = "'3.541.21"
$ws.Range('E8').Value = "'  +0.99%  "
$ws.Range('E10').Value = "'  +5.75%  "
$ws.Range('D11').Value = "'6.75"
$ws.Range('E11').Value = "'  -0.18%  "
$ws.Range('D12').Value = "'0.588"
$ws.Range('E12').Value = "'  -1.19%  "
$ws.Range('D13').Value = "'47.83"
$ws.Range('E13').Value = "'  +1.74%  "
$ws.Range('E14').Value = "'  +1.60%  "
$ws.Range('D15').Value = "'4.120.03"
$ws.Range('E15').Value = "'  +1.08%  "
$ws.Range('D16').Value = "'630.70"
$ws.Range('E16').Value = "'  -6.93%  "
$ws.Range('D17').Value = "'8.46"
$ws.Range('E17').Value = "'  -3.22%  "
$ws.Range('D18').Value = "'70.810.67"
$ws.Range('E18').Value = "'  +2.23%  "
$ws.Range('D19').Value = "'3.555.56"
$ws.Range('E19').Value = "'  +1.20%  "
$ws.Range('E20').Value = "'  -1.80%  "
$ws.Range('E21').Value = "'  +0.22%  "
$ws.Range('D22').Value = "'10.10"
$ws.Range('E22').Value = "'  -9.92%  "
$ws.Range('D23').Value = "'0.891"
$ws.Range('E23').Value = "'  -1.42%  "
$ws.Range('D24').Value = "'15.99"
$ws.Range('E24').Value = "'  -0.88%  "
$ws.Range('D25').Value = "'97.26"
$ws.Range('E25').Value = "'  -0.78%  "
$ws.Range('E26').Value = "'  -0.09%  "
$ws.Range('E27').Value = "'  -0.10%  "
$ws.Range('D28').Value = "'2.63"
$ws.Range('E28').Value = "'  -0.98%  "
$ws.Range('D29').Value = "'9.29"
$ws.Range('E29').Value = "'  -1.70%  "
$ws.Range('D30').Value = "'33.51"
$ws.Range('E30').Value = "'  +1.46%  "
$ws.Range('E31').Value = "'  -1.93%  "
$ws.Range('D32').Value = "'8.50"
$ws.Range('E32').Value = "'  -2.66%  "
$ws.Range('E33').Value = "'  -0.88%  "
$ws.Range('E34').Value = "'  -2.88%  "
$ws.Range('D35').Value = "'570.73"
$ws.Range('E35').Value = "'  -4.17%  "
$ws.Range('D36').Value = "'3.67"
$ws.Range('E36').Value = "'  +1.50%  "
$ws.Range('D37').Value = "'10.82"
$ws.Range('E37').Value = "'  -0.82%  "
$ws.Range('E38').Value = "'  -1.82%  "
$ws.Range('D39').Value = "'57.52"
$ws.Range('E39').Value = "'  +0.38%  "
$ws.Range('D40').Value = "'0.999"
$ws.Range('E40').Value = "'  -0.01%  "
$ws.Range('D41').Value = "'0.143"
$ws.Range('E41').Value = "'  +5.47%  "
$ws.Range('D42').Value = "'0.0455"
$ws.Range('E42').Value = "'  +3.66%  "
$ws.Range('E43').Value = "'  -1.50%  "
$ws.Range('D44').Value = "'3.350.84"
$ws.Range('E44').Value = "'  -1.88%  "
$ws.Range('E45').Value = "'  +5.22%  "
$ws.Range('D46').Value = "'0.0₃0723"
$ws.Range('E46').Value = "'  +1.97%  "
$ws.Range('D47').Value = "'33.29"
$ws.Range('E47').Value = "'  -0.25%  "
$ws.Range('E48').Value = "'  +2.79%  "
$ws.Range('E49').Value = "'  -2.13%  "
$ws.Range('D50').Value = "'134.96"
$ws.Range('E50').Value = "'  +2.15%  "
$ws.Range('E51').Value = "'  -1.86%  "
